# Add files via upload
# Populate columns D (Info to display) and E (Numero de serie) on the
# "Liste" sheet, and update a couple of existing "Is Available" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Liste")

# Update a couple of existing "Is Available" cells
$ws.Range("C3").Value = "N/a"
$ws.Range("C4").Value = "Macron démission"

# New "Info to display" column (header + values, D2 filled in last)
$ws.Range("D1").Value = "Info to display"
$ws.Range("D3").Value = "a = 3"
$ws.Range("D4").Value = "a = 4"
$ws.Range("D5").Value = "a = 5"
$ws.Range("D6").Value = "a = 6"
$ws.Range("D2").Value = "a = 2"

# New "Numero de serie" column
$ws.Range("E1").Value = "Numero de serie"
$ws.Range("E2").Value = 21413
$ws.Range("E3").Value = 15345
$ws.Range("E4").Value = 5216534
$ws.Range("E5").Value = 153
$ws.Range("E6").Value = 54848

# Selection moves to E7 after the edits
$ws.Range("E7").Select() | Out-Null
